$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing last header cell (G1) onto the new
# header cell (H1) so the new "Save" column header matches the other
# header cells (bold font, border, centered/top aligned).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new header text and the corresponding data value.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
